# Update "paises.xlsx" (COVID data) — refresh timestamp and per-country statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 03:09"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1263092
$ws.Range("C4").Value = 25459
$ws.Range("D4").Value = 212981
$ws.Range("E4").Value = 975312
$ws.Range("F4").Value = 15827
$ws.Range("G4").Value = 2528
$ws.Range("H4").Value = 74799

# Row 12 - Brasil
$ws.Range("B12").Value = 126611
$ws.Range("C12").Value = 11896
$ws.Range("D12").Value = 51370
$ws.Range("E12").Value = 66653
$ws.Range("F12").Value = 8318
$ws.Range("G12").Value = 667
$ws.Range("H12").Value = 8588

# Row 15 - Canada
$ws.Range("D15").Value = 28171
$ws.Range("E15").Value = 31093

# Row 47 - Noruega
$ws.Range("B47").Value = 7996
$ws.Range("C47").Value = 41
$ws.Range("E47").Value = 7748

# Row 66 - Tailandia
$ws.Range("B66").Value = 3091
$ws.Range("C66").Value = 372
$ws.Range("D66").Value = 294
$ws.Range("E66").Value = 2779
$ws.Range("F66").Value = 4
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 18

# Row 67 - Oman
$ws.Range("B67").Value = 2989
$ws.Range("C67").Value = 1
$ws.Range("D67").Value = 2761
$ws.Range("E67").Value = 173
$ws.Range("F67").Value = 61
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 55

# Row 68 - Armenia
$ws.Range("B68").Value = 2903
$ws.Range("C68").Value = 168
$ws.Range("D68").Value = 888
$ws.Range("E68").Value = 2002
$ws.Range("F68").Value = 17
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 13

# Row 69 - Ghana
$ws.Range("B69").Value = 2782
$ws.Range("C69").Value = 163
$ws.Range("D69").Value = 1135
$ws.Range("E69").Value = 1607
$ws.Range("F69").Value = 10
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 40

# Row 128 - Estado de Palestina
$ws.Range("B128").Value = 379
$ws.Range("C128").Value = 18
$ws.Range("D128").Value = 176
$ws.Range("E128").Value = 193
$ws.Range("F128").Value = 1
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 10

# Row 129 - Venezuela
$ws.Range("B129").Value = 374
$ws.Range("C129").Value = 3
$ws.Range("D129").Value = 174
$ws.Range("E129").Value = 198
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 2
